# "adloori to davuluri completed"
# Fill in the "Points for grading" (column E) scores for the first two
# rubric sections (Generic / Customer Class) to match the "Total Points"
# (column D) already entered for those rows, then leave the selection on
# E15 (the subtotal for the section just completed) as the grader did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Generic section (rows 3-6)
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2

# Customer Class section (rows 10-14)
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2

# Move the active selection to E15, matching where grading left off.
$ws.Range("E15").Select()
